$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "27.196.54"
$ws.Range("E2").Value = "  +1.30%  "
$ws.Range("D3").Value = "1.643.16"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "217.24"
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("E6").Value = "  +1.36%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E9").Value = "  +1.02%  "
$ws.Range("D10").Value = "19.93"
$ws.Range("E10").Value = "  +0.44%  "
$ws.Range("D11").Value = "0.0847"
$ws.Range("E11").Value = "  +0.37%  "
$ws.Range("D12").Value = "1.873.51"
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("D13").Value = "1.671.94"
$ws.Range("E13").Value = "  +1.75%  "
$ws.Range("D14").Value = "4.14"
$ws.Range("E14").Value = "  +0.44%  "
$ws.Range("E15").Value = "  +2.95%  "
$ws.Range("D16").Value = "67.36"
$ws.Range("E16").Value = "  +1.53%  "
$ws.Range("D17").Value = "27.196.28"
$ws.Range("E17").Value = "  +1.21%  "
$ws.Range("E18").Value = "  +1.03%  "
$ws.Range("D19").Value = "218.88"
$ws.Range("E19").Value = "  +0.59%  "
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").Value = "6.81"
$ws.Range("E21").Value = "  +2.96%  "
$ws.Range("D22").Value = "2.56"
$ws.Range("E22").Value = "  +4.93%  "
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").Value = "147.66"
$ws.Range("E25").Value = "  +1.56%  "
$ws.Range("D26").Value = "7.54"
$ws.Range("E26").Value = "  +2.44%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").Value = "15.74"
$ws.Range("E29").Value = "  -0.61%  "
$ws.Range("D31").Value = "1.18"
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("D32").Value = "3.39"
$ws.Range("E32").Value = "  +0.65%  "
$ws.Range("D33").Value = "3.02"
$ws.Range("E33").Value = "  +1.03%  "
$ws.Range("E34").Value = "  +0.73%  "
$ws.Range("D35").Value = "1.262.01"
$ws.Range("E35").Value = "  +1.43%  "
$ws.Range("E36").Value = "  +0.92%  "
$ws.Range("E37").Value = "  +2.08%  "
$ws.Range("E38").Value = "  +1.30%  "
$ws.Range("D39").Value = "0.849"
$ws.Range("E39").Value = "  +1.78%  "
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("D41").Value = "0.807"
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("E42").Value = "  +6.22%  "
$ws.Range("D43").Value = "5.36"
$ws.Range("E43").Value = "  -0.15%  "
$ws.Range("D44").Value = "1.784.43"
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("D45").Value = "61.72"
$ws.Range("E45").Value = "  +1.23%  "
$ws.Range("D46").Value = "91.61"
$ws.Range("E46").Value = "  +0.20%  "
$ws.Range("D47").Value = "1.60"
$ws.Range("E47").Value = "  +1.01%  "
$ws.Range("D48").Value = "0.0₆0105"
$ws.Range("E48").Value = "  -0.49%  "
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("D50").Value = "7.65"
$ws.Range("E50").Value = "  +1.43%  "
$ws.Range("E51").Value = "  +0.32%  "
